$wb = $excel.ActiveWorkbook

# "Building upgrade" cost sheets: column D header was mislabeled with the
# FLOAT_* key names (speedup/resource) even though the stored values are
# plain integers. Rename the headers to the correct INT_* keys.

$ws = $wb.Worksheets.Item("time")
$ws.Activate()
$ws.Range("D1").Value = "INT_speedup"
[void]$ws.Range("D2").Select()

$ws = $wb.Worksheets.Item("wood")
$ws.Activate()
$ws.Range("D1").Value = "INT_resource"
[void]$ws.Range("D2").Select()

$ws = $wb.Worksheets.Item("stone")
$ws.Activate()
$ws.Range("D1").Value = "INT_resource"
[void]$ws.Range("D1").Select()

$ws = $wb.Worksheets.Item("iron")
$ws.Activate()
$ws.Range("D1").Value = "INT_resource"
[void]$ws.Range("D1").Select()

$ws = $wb.Worksheets.Item("food")
$ws.Activate()
$ws.Range("D1").Value = "INT_resource"
[void]$ws.Range("D1").Select()

$ws = $wb.Worksheets.Item("coin")
$ws.Activate()
$ws.Range("D1").Value = "INT_resource"
[void]$ws.Range("D1").Select()

# Restore the original active sheet
$ws = $wb.Worksheets.Item("material")
$ws.Activate()
